$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("County")

$rows = @(65, 66, 67)
foreach ($r in $rows) {
    $cellB = $ws.Range("B$r")
    $cellB.NumberFormat = "@"
    $cellB.Value = "0.00%"
    $cellB.Style = "Normal"

    $cellC = $ws.Range("C$r")
    $cellC.NumberFormat = "@"
    $cellC.Value = "`$0"
    $cellC.Style = "Normal"

    $cellD = $ws.Range("D$r")
    $cellD.NumberFormat = "@"
    $cellD.Value = "0.00%"
    $cellD.Style = "Normal"

    $cellE = $ws.Range("E$r")
    $cellE.NumberFormat = "@"
    $cellE.Value = "0.00%"
    $cellE.Style = "Normal"

    $cellF = $ws.Range("F$r")
    $cellF.NumberFormat = "@"
    $cellF.Value = "0.00%"
    $cellF.Style = "Normal"
}
